$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: Oct 10, 2024, 18:30 - 18:58
$ws.Range("A10").Value = 45575
$ws.Range("B10").Value = 0.77083333333333337
$ws.Range("C10").Value = 0.79027777777777775

# Row 11: Oct 10, 2024, 20:00 - 21:00
$ws.Range("A11").Value = 45575
$ws.Range("B11").Value = 0.83333333333333337
$ws.Range("C11").Value = 0.875

# Row 4: add M4, N4 formulas
$ws.Range("M4").Formula = "=SUM(D10:D11)"
$ws.Range("N4").Formula = "=SUM(G10:G11)"

# Update selection to I18 (as seen in diff)
$ws.Range("I18").Select()
